$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row above the current row 5 ("Number of exp. conditions...")
# to hold the new "LP solver" setting, shifting everything below down by one row.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "LP solver (linprog or gurobi)"
$ws.Range("B5").Value = "gurobi"

# Match the header-label style used in column A (bold border style) but
# left-aligned instead of centered, like the other "value" rows that use
# left alignment for string answers.
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4160

# Restore selection/active-cell state to point at the newly inserted row,
# and make the general sheet the active tab (mirrors the workbook view
# change recorded in the edit).
$ws.Activate()
$ws.Range("A5:B5").Select()
